$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 14580.2
$ws.Range("J69").Value = 13225.25
$ws.Range("L69").Value = 39675.75
$ws.Range("N69").Value = -41423.75
$ws.Range("H72").Value = 14580.2
$ws.Range("J72").Value = 13225.25
$ws.Range("L72").Value = 119027.25
$ws.Range("N72").Value = -127763.25
$ws.Range("H80").Value = 1576.4615
$ws.Range("I80").Value = 1332.8889
$ws.Range("K80").Value = 3998.6667
$ws.Range("M80").Value = -3000.6667
$ws.Range("H83").Value = 1576.4615
$ws.Range("I83").Value = 1332.8889
$ws.Range("K83").Value = 11996.0001
$ws.Range("M83").Value = -7004.000099999999
$ws.Range("I86").Value = 1619.5
$ws.Range("K86").Value = 1619.5
$ws.Range("M86").Value = -496.5
$ws.Range("I89").Value = 1619.5
$ws.Range("K89").Value = 8097.5
$ws.Range("M89").Value = -2481.5
$ws.Range("H100").Value = 7142.706
$ws.Range("I100").Value = 2643.4
$ws.Range("K100").Value = 2643.4
$ws.Range("M100").Value = -2102.4
$ws.Range("H112").Value = 1569.4814
$ws.Range("J112").Value = 1592.1923
$ws.Range("L112").Value = 4776.5769
$ws.Range("N112").Value = -6992.5769
$ws.Range("H137").Value = 3757.6086
$ws.Range("I137").Value = 2354.4
$ws.Range("K137").Value = 7063.200000000001
$ws.Range("M137").Value = -4513.200000000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H74").Value = 1427.4193
$ws.Range("I74").Value = 937.7143
$ws.Range("K74").Value = 937.7143
$ws.Range("M74").Value = -63.71429999999998
$ws.Range("H77").Value = 1427.4193
$ws.Range("I77").Value = 937.7143
$ws.Range("K77").Value = 4688.5715
$ws.Range("M77").Value = -320.5715
$ws.Range("H102").Value = 2291.2222
$ws.Range("I102").Value = 2131.8823
$ws.Range("K102").Value = 2131.8823
$ws.Range("M102").Value = -509.8823000000002

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 1001.75
$ws.Range("I7").Value = 669
$ws.Range("K7").Value = 669
$ws.Range("M7").Value = -556
$ws.Range("H60").Value = 52444.75
$ws.Range("J60").Value = 52444.75
$ws.Range("L60").Value = 52444.75
$ws.Range("N60").Value = -53642.75
$ws.Range("H86").Value = 2495.6924
$ws.Range("I86").Value = 2273.7896
$ws.Range("K86").Value = 2273.7896
$ws.Range("M86").Value = -1150.7896
$ws.Range("H89").Value = 2495.6924
$ws.Range("I89").Value = 2273.7896
$ws.Range("K89").Value = 11368.948
$ws.Range("M89").Value = -5752.948
$ws.Range("H102").Value = 5547
$ws.Range("I102").Value = 5547
$ws.Range("K102").Value = 5547
$ws.Range("M102").Value = -2302

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 4748.5
$ws.Range("I3").Value = 3995
$ws.Range("K3").Value = 3995
$ws.Range("M3").Value = -3882
$ws.Range("H20").Value = 70471.5
$ws.Range("J20").Value = 70471.5
$ws.Range("L20").Value = 70471.5
$ws.Range("N20").Value = -70943.5
$ws.Range("H30").Value = 70471.5
$ws.Range("J30").Value = 70471.5
$ws.Range("L30").Value = 70471.5
$ws.Range("N30").Value = -70653.5
$ws.Range("H31").Value = 2989.1875
$ws.Range("I31").Value = 1576.3636
$ws.Range("J31").Value = 6097.4
$ws.Range("K31").Value = 1576.3636
$ws.Range("L31").Value = 6097.4
$ws.Range("M31").Value = -1281.3636
$ws.Range("N31").Value = -6687.4
$ws.Range("H34").Value = 2989.1875
$ws.Range("I34").Value = 1576.3636
$ws.Range("J34").Value = 6097.4
$ws.Range("K34").Value = 1576.3636
$ws.Range("L34").Value = 6097.4
$ws.Range("M34").Value = -1374.3636
$ws.Range("N34").Value = -6501.4
$ws.Range("H105").Value = 1537.95
$ws.Range("I105").Value = 1476.3334
$ws.Range("K105").Value = 1476.3334
$ws.Range("M105").Value = 270.6666
$ws.Range("H128").Value = 70471.5
$ws.Range("J128").Value = 70471.5
$ws.Range("L128").Value = 70471.5
$ws.Range("N128").Value = -80431.5
$ws.Range("H134").Value = 3878.7896
$ws.Range("I134").Value = 3420.077
$ws.Range("J134").Value = 4872.6665
$ws.Range("K134").Value = 10260.231
$ws.Range("L134").Value = 14617.9995
$ws.Range("M134").Value = -7725.231
$ws.Range("N134").Value = -19687.9995

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 45460656
$ws.Range("I70").Value = 4666.25
$ws.Range("J70").Value = 55561988
$ws.Range("K70").Value = 4666.25
$ws.Range("L70").Value = 55561988
$ws.Range("M70").Value = -4396.25
$ws.Range("N70").Value = -55562528
$ws.Range("H73").Value = 45460656
$ws.Range("I73").Value = 4666.25
$ws.Range("J73").Value = 55561988
$ws.Range("K73").Value = 4666.25
$ws.Range("L73").Value = 55561988
$ws.Range("M73").Value = -3730.25
$ws.Range("N73").Value = -55563860
$ws.Range("H113").Value = 395311.8
$ws.Range("J113").Value = 12447.5
$ws.Range("L113").Value = 12447.5
$ws.Range("N113").Value = -16787.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 588.3158
$ws.Range("I55").Value = 504.81818
$ws.Range("J55").Value = 703.125
$ws.Range("K55").Value = 504.81818
$ws.Range("L55").Value = 703.125
$ws.Range("M55").Value = -331.81818
$ws.Range("N55").Value = -1049.125
$ws.Range("H100").Value = 6520.2354
$ws.Range("I100").Value = 2523.0667
$ws.Range("K100").Value = 2523.0667
$ws.Range("M100").Value = -1982.0667
$ws.Range("H115").Value = 68975
$ws.Range("J115").Value = 68975
$ws.Range("L115").Value = 68975
$ws.Range("N115").Value = -71325
$ws.Range("H137").Value = 75000
$ws.Range("J137").Value = 75000
$ws.Range("L137").Value = 75000
$ws.Range("N137").Value = -85200

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 11724.833
$ws.Range("I41").Value = 11571
$ws.Range("K41").Value = 11571
$ws.Range("M41").Value = -11181
$ws.Range("H126").Value = 3053.6538
$ws.Range("I126").Value = 3035.7856
$ws.Range("J126").Value = 3074.5
$ws.Range("K126").Value = 9107.356800000001
$ws.Range("L126").Value = 9223.5
$ws.Range("M126").Value = -6637.356800000001
$ws.Range("N126").Value = -14163.5
$ws.Range("H132").Value = 3458.353
$ws.Range("I132").Value = 2372.0454
$ws.Range("K132").Value = 7116.1362
$ws.Range("M132").Value = -4586.1362
